# 6MWT collection sheet: insert a new "categories" column (column K) ahead of
# "short_name", shifting every subsequent column one to the right, then keep
# the AutoFilter / filter-database defined name / selection in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K (pushes existing K:AG -> L:AH).
$ws.Columns("K").Insert()

# New header text for the inserted column.
$ws.Range("K1").Value2 = "categories"

# Give the new column the same custom width (~13 chars) the other
# non-bestFit data columns (e.g. G) use, instead of the sheet default.
$ws.Columns("K").ColumnWidth = 12.17

# The AutoFilter range doesn't auto-grow with the column insert in this
# engine, so turn it off and reapply across the new, wider range.
$ws.AutoFilterMode = $false
$ws.Range("A1:AH17").AutoFilter()

# Keep the hidden _FilterDatabase defined name (driven by the AutoFilter) in
# sync with the new range too.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Collection_QRS_6MWT!_FilterDatabase") {
        $n.RefersTo = "=Collection_QRS_6MWT!`$A`$1:`$AH`$17"
    }
}

# Match the author's final selection/scroll position on the sheet.
$ws.Range("J25").Select()
